# Move all non-dispatchable resources (geothermal, biomass, municipal solid
# waste, all other renewables) to guaranteed dispatch on the BGDPbES sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGDPbES")

# biomass (row 9) -> guaranteed dispatch = 1
$ws.Range("B9").Value = 1

# geothermal (row 10) -> guaranteed dispatch = 1
$ws.Range("B10").Value = 1

# municipal solid waste (row 17) -> guaranteed dispatch = 1, and this row
# previously held static zeros with no formulas; give it the same
# "carry forward" formula pattern used by the other rows.
$ws.Range("B17").Value = 1
$ws.Range("C17").Formula = "=B17"
$ws.Range("D17:AK17").Formula = "=C17"

# Leave the cursor on B6 like the saved workbook did.
[void]$ws.Range("B6").Select()
